$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.032.79'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.940.76'
$ws.Range("E3").Value = '  +1.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '377.10'
$ws.Range("E5").Value = '  +2.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.15'
$ws.Range("E6").Value = '  +1.55%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.19%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.73'
$ws.Range("E10").Value = '  -0.14%  '

$ws.Range("E11").Value = '  +0.39%  '

$ws.Range("E12").Value = '  +0.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.402.04'
$ws.Range("E13").Value = '  +1.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.27'
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.43'
$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.946.54'
$ws.Range("E16").Value = '  +1.87%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.951'
$ws.Range("E17").Value = '  +1.32%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.966.01'
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.32'
$ws.Range("E19").Value = '  +1.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.30'
$ws.Range("E20").Value = '  +1.04%  '

$ws.Range("E21").Value = '  -0.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0955'
$ws.Range("E22").Value = '  +1.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.82'
$ws.Range("E23").Value = '  +0.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.00'
$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("E25").Value = '  +3.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.168'
$ws.Range("E26").Value = '  -0.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.14'
$ws.Range("E27").Value = '  +18.40%  '

$ws.Range("E28").Value = '  +0.05%  '

$ws.Range("E29").Value = '  +3.52%  '

$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("E31").Value = '  +6.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.77'
$ws.Range("E32").Value = '  -1.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.39'
$ws.Range("E33").Value = '  -1.50%  '

$ws.Range("E34").Value = '  -2.29%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.85'
$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0444'
$ws.Range("E36").Value = '  +6.31%  '

$ws.Range("E37").Value = '  +0.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.05'
$ws.Range("E38").Value = '  -0.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.10'
$ws.Range("E39").Value = '  +1.11%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.57'
$ws.Range("E40").Value = '  -3.42%  '

$ws.Range("E41").Value = '  +2.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.83'
$ws.Range("E42").Value = '  -0.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '121.99'
$ws.Range("E43").Value = '  +3.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.78'
$ws.Range("E44").Value = '  -1.21%  '

$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.07'
$ws.Range("E45").Value = '  -1.33%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.280'
$ws.Range("E46").Value = '  +16.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.38'
$ws.Range("E47").Value = '  +2.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.031.26'
$ws.Range("E48").Value = '  -0.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.19'
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0343'
$ws.Range("E50").Value = '  +9.95%  '

$ws.Range("E51").Value = '  +0.57%  '
